$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.000.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -6.37%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.278.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.62%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "558.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.57%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "127.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.98%  "

# Row 7
$ws.Range("E7").Value = "  +0.05%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.278.88"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.61%  "

# Row 9
$ws.Range("E9").Value = "  -2.18%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.84%  "

# Row 11
$ws.Range("E11").Value = "  -5.53%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.368"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.78%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.842.50"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.52%  "

# Row 14
$ws.Range("E14").Value = "  -0.30%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.291.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.32%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000166"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.48%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.279.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -5.91%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "24.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.34%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.40%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.24%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -10.56%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "349.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.550"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.08%  "

# Row 24
$ws.Range("E24").Value = "  -0.27%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.414.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.52%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.52%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000107"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.79%  "

# Row 28
$ws.Range("E28").Value = "  +0.26%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.18"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.95%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.39%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.23%  "

# Row 32
$ws.Range("E32").Value = "  -6.39%  "

# Row 33
$ws.Range("E33").Value = "  +0.04%  "

# Row 34
$ws.Range("E34").Value = "  -2.42%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.309.82"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.50%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "22.62"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.36%  "

# Row 37
$ws.Range("E37").Value = "  -1.13%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.71"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.88%  "

# Row 39
$ws.Range("E39").Value = "  -2.52%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "157.11"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.07%  "

# Row 41
$ws.Range("E41").Value = "  -4.51%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.48%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.32"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.40%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.738"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -7.40%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.25%  "

# Row 47
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.16%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -5.31%  "

# Row 49
$ws.Range("E49").Value = "  -1.10%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.854"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -5.38%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.58%  "

